# Update the division-fact answer table in place.
# The table has 20 rows total: 5 "data" rows (1, 5, 9, 13, 17) each
# holding 5 worked-division answers, interleaved with 3 blank rows.
# Each data cell's text is replaced with its new answer while leaving
# every other part of the document (fonts, paragraph/table structure,
# the date heading, and the blank rows) untouched.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 1
$table.Cell(1, 1).Range.Text = "19÷8=2, 3"
$table.Cell(1, 2).Range.Text = "46÷9=5, 1"
$table.Cell(1, 3).Range.Text = "80÷6=13, 2"
$table.Cell(1, 4).Range.Text = "33÷2=16, 1"
$table.Cell(1, 5).Range.Text = "64÷6=10, 4"

# Row 5
$table.Cell(5, 1).Range.Text = "66÷2=33, 0"
$table.Cell(5, 2).Range.Text = "77÷8=9, 5"
$table.Cell(5, 3).Range.Text = "61÷4=15, 1"
$table.Cell(5, 4).Range.Text = "61÷4=15, 1"
$table.Cell(5, 5).Range.Text = "63÷5=12, 3"

# Row 9
$table.Cell(9, 1).Range.Text = "35÷9=3, 8"
$table.Cell(9, 2).Range.Text = "13÷8=1, 5"
$table.Cell(9, 3).Range.Text = "32÷9=3, 5"
$table.Cell(9, 4).Range.Text = "29÷8=3, 5"
$table.Cell(9, 5).Range.Text = "30÷4=7, 2"

# Row 13
$table.Cell(13, 1).Range.Text = "84÷6=14, 0"
$table.Cell(13, 2).Range.Text = "12÷7=1, 5"
$table.Cell(13, 3).Range.Text = "20÷4=5, 0"
$table.Cell(13, 4).Range.Text = "31÷3=10, 1"
$table.Cell(13, 5).Range.Text = "83÷3=27, 2"

# Row 17
$table.Cell(17, 1).Range.Text = "56÷2=28, 0"
$table.Cell(17, 2).Range.Text = "87÷4=21, 3"
$table.Cell(17, 3).Range.Text = "88÷2=44, 0"
$table.Cell(17, 4).Range.Text = "98÷2=49, 0"
$table.Cell(17, 5).Range.Text = "83÷8=10, 3"

